$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1, J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from an existing header cell (H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-28
$values = @{
    2  = @(6, 6)
    3  = @(7, 7)
    4  = @(5, 5)
    5  = @(7, 7)
    6  = @(5, 7)
    7  = @(6, 7)
    8  = @(9, 10)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(7, 7)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(8, 8)
    19 = @(5, 6)
    20 = @(5, 5)
    21 = @(7, 7)
    22 = @(6, 6)
    23 = @(9, 9)
    24 = @(4, 4)
    25 = @(5, 5)
    26 = @(6, 6)
    27 = @(5, 5)
    28 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
